$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Roxie" harmonic values (columns R:AF) for rows 2,3,4,6,7,8
# Commit: Scan profiles plot changed from Bn to bn

# Row 2
$ws.Cells.Item(2, 18).Value = 0.01463892547037577
$ws.Cells.Item(2, 19).Value = 0.00141215296033211
$ws.Cells.Item(2, 20).Value = 0.00005454198232676336
$ws.Cells.Item(2, 21).Value = -0.0000008040142526557588
$ws.Cells.Item(2, 22).Value = 0.0000009830035190006227
$ws.Cells.Item(2, 23).Value = -0.0000009245954587075519
$ws.Cells.Item(2, 24).Value = 0.0000007419211018965311
$ws.Cells.Item(2, 25).Value = 0.0000006643875594413854
$ws.Cells.Item(2, 26).Value = 0.0000009993268171677071
$ws.Cells.Item(2, 27).Value = -0.000002542159181247002
$ws.Cells.Item(2, 28).Value = -0.0000001063000631743917
$ws.Cells.Item(2, 29).Value = -0.0000006277928083875839
$ws.Cells.Item(2, 30).Value = 0.000001938132214561294
$ws.Cells.Item(2, 31).Value = -0.0000003946715831337595
$ws.Cells.Item(2, 32).Value = 0.000002339967407592905

# Row 3
$ws.Cells.Item(3, 18).Value = 0.1915074386132086
$ws.Cells.Item(3, 19).Value = 0.007564351987168769
$ws.Cells.Item(3, 20).Value = 0.01510192915113025
$ws.Cells.Item(3, 21).Value = 0.00001931235276659419
$ws.Cells.Item(3, 22).Value = -0.0008751033329785377
$ws.Cells.Item(3, 23).Value = -0.00003064148130251899
$ws.Cells.Item(3, 24).Value = -0.0001265267828990458
$ws.Cells.Item(3, 25).Value = -0.000001000858501059038
$ws.Cells.Item(3, 26).Value = 0.00002597608896544229
$ws.Cells.Item(3, 27).Value = -0.0000009064687956694585
$ws.Cells.Item(3, 28).Value = 0.000001643409258886945
$ws.Cells.Item(3, 29).Value = -0.00000006971314100560452
$ws.Cells.Item(3, 30).Value = -0.00001657249320790441
$ws.Cells.Item(3, 31).Value = 0.000002093482330864965
$ws.Cells.Item(3, 32).Value = -0.0000189016835541236

# Row 4
$ws.Cells.Item(4, 18).Value = 0.4493927569752659
$ws.Cells.Item(4, 19).Value = 0.01672010369830522
$ws.Cells.Item(4, 20).Value = 0.04170324148176154
$ws.Cells.Item(4, 21).Value = -0.0008127587750616519
$ws.Cells.Item(4, 22).Value = -0.002473859915468567
$ws.Cells.Item(4, 23).Value = -0.0001685251631970967
$ws.Cells.Item(4, 24).Value = -0.0004027237224865842
$ws.Cells.Item(4, 25).Value = -0.000001339009105704611
$ws.Cells.Item(4, 26).Value = 0.00001312882368526758
$ws.Cells.Item(4, 27).Value = 0.000006630778455479374
$ws.Cells.Item(4, 28).Value = 0.000007225681601377494
$ws.Cells.Item(4, 29).Value = -0.0000032990976709217
$ws.Cells.Item(4, 30).Value = 0.000004979468139097218
$ws.Cells.Item(4, 31).Value = -0.000007409051268352657
$ws.Cells.Item(4, 32).Value = -0.00000004687450102856368

# Row 6
$ws.Cells.Item(6, 18).Value = 0.7336463186779769
$ws.Cells.Item(6, 19).Value = 0.02070657828966237
$ws.Cells.Item(6, 20).Value = 0.05192458508802378
$ws.Cells.Item(6, 21).Value = -0.0006427893663666269
$ws.Cells.Item(6, 22).Value = -0.006063730792680811
$ws.Cells.Item(6, 23).Value = -0.0001669002519611993
$ws.Cells.Item(6, 24).Value = -0.00323903921839681
$ws.Cells.Item(6, 25).Value = -0.000003557158391825837
$ws.Cells.Item(6, 26).Value = 0.002933944641203898
$ws.Cells.Item(6, 27).Value = 0.00000598095564564697
$ws.Cells.Item(6, 28).Value = -0.001585275761174224
$ws.Cells.Item(6, 29).Value = -0.000003449716212711299
$ws.Cells.Item(6, 30).Value = 0.000612774247380056
$ws.Cells.Item(6, 31).Value = -0.000003592598394549926
$ws.Cells.Item(6, 32).Value = -0.0001672795585071591

# Row 7
$ws.Cells.Item(7, 18).Value = 1.448065781083953
$ws.Cells.Item(7, 19).Value = 0.005571209329262463
$ws.Cells.Item(7, 20).Value = 1.163443047236701
$ws.Cells.Item(7, 21).Value = 0.00009159232927933242
$ws.Cells.Item(7, 22).Value = -0.5765314700226601
$ws.Cells.Item(7, 23).Value = -0.00007542552327761188
$ws.Cells.Item(7, 24).Value = -0.1126368539820084
$ws.Cells.Item(7, 25).Value = 0.00000268586666216847
$ws.Cells.Item(7, 26).Value = 0.07243936636103107
$ws.Cells.Item(7, 27).Value = -0.000002737789357151291
$ws.Cells.Item(7, 28).Value = -0.08773909671812553
$ws.Cells.Item(7, 29).Value = 0.0000006210903687960008
$ws.Cells.Item(7, 30).Value = 0.04300646116931833
$ws.Cells.Item(7, 31).Value = 0.000001060236444961737
$ws.Cells.Item(7, 32).Value = 0.004177124478312502

# Row 8
$ws.Cells.Item(8, 18).Value = 0.0483354075648443
$ws.Cells.Item(8, 19).Value = 0.0009566595609196547
$ws.Cells.Item(8, 20).Value = -0.0001731044517259024
$ws.Cells.Item(8, 21).Value = 0.0000001290428993134621
$ws.Cells.Item(8, 22).Value = -0.000002753215289616471
$ws.Cells.Item(8, 23).Value = -0.000002601357981481483
$ws.Cells.Item(8, 24).Value = -0.000002279049546161856
$ws.Cells.Item(8, 25).Value = 0.00000009416919247788845
$ws.Cells.Item(8, 26).Value = 0.0000004805505696602163
$ws.Cells.Item(8, 27).Value = 0.000002030070236388733
$ws.Cells.Item(8, 28).Value = -0.000004621547622708025
$ws.Cells.Item(8, 29).Value = -0.000001305868100098961
$ws.Cells.Item(8, 30).Value = -0.000005368057707662112
$ws.Cells.Item(8, 31).Value = -0.000002176308543964261
$ws.Cells.Item(8, 32).Value = 0.000001786583899709863
